$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.630.21"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").Value = "3.633.17"
$ws.Range("E3").Value = "  +1.07%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +14.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "653.46"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.417"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.40%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.07"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.02%  "

$ws.Range("D11").Value = "3.630.87"
$ws.Range("E11").Value = "  +1.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.20%  "

$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.80%  "

$ws.Range("D15").Value = "4.307.20"
$ws.Range("E15").Value = "  +1.00%  "

$ws.Range("D16").Value = "96.414.42"
$ws.Range("E16").Value = "  -0.17%  "

$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +11.72%  "

$ws.Range("D19").Value = "3.628.10"
$ws.Range("E19").Value = "  +1.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.527"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.77%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "510.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.75%  "

$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.12%  "

$ws.Range("E25").Value = "  +1.96%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.96%  "

$ws.Range("E29").Value = "  +13.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.53%  "

$ws.Range("E32").Value = "  -0.12%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.70%  "

$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.21%  "

$ws.Range("E36").Value = "  +8.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.578"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.16%  "

$ws.Range("E38").Value = "  +3.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "612.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.20%  "

$ws.Range("E40").Value = "  +2.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.954"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.57%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.33%  "

$ws.Range("B43").Value = "ImmutableX"
$ws.Range("C43").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.21%  "

$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.38%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0442"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.45%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.27%  "

$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.17%  "

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.408"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +15.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.22%  "

